$wb = $excel.ActiveWorkbook

# --- Sheet ALC (index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H42").Value = [double]3456.6667
$ws.Range("I42").Value = [double]1584
$ws.Range("J42").Value = [double]5329.3335
$ws.Range("K42").Value = [double]4752
$ws.Range("L42").Value = [double]15988.0005
$ws.Range("M42").Value = [double]-4522
$ws.Range("N42").Value = [double]-16448.0005
$ws.Range("H107").Value = [double]809.875
$ws.Range("I107").Value = [double]695.5
$ws.Range("K107").Value = [double]695.5
$ws.Range("M107").Value = [double]1224.5
$ws.Range("H118").Value = [double]1093.4
$ws.Range("I118").Value = [double]1093.4
$ws.Range("K118").Value = [double]3280.2
$ws.Range("M118").Value = [double]-1623.2
$ws.Range("H125").ClearContents()
$ws.Range("I125").ClearContents()
$ws.Range("J125").ClearContents()
$ws.Range("K125").ClearContents()
$ws.Range("L125").ClearContents()
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()
$ws.Range("H126").ClearContents()
$ws.Range("I126").ClearContents()
$ws.Range("J126").ClearContents()
$ws.Range("K126").ClearContents()
$ws.Range("L126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H127").ClearContents()
$ws.Range("I127").ClearContents()
$ws.Range("J127").ClearContents()
$ws.Range("K127").ClearContents()
$ws.Range("L127").ClearContents()
$ws.Range("M127").ClearContents()
$ws.Range("H128").ClearContents()
$ws.Range("I128").ClearContents()
$ws.Range("J128").ClearContents()
$ws.Range("K128").ClearContents()
$ws.Range("L128").ClearContents()
$ws.Range("H129").ClearContents()
$ws.Range("I129").ClearContents()
$ws.Range("J129").ClearContents()
$ws.Range("K129").ClearContents()
$ws.Range("L129").ClearContents()
$ws.Range("M129").ClearContents()
$ws.Range("N129").ClearContents()
$ws.Range("H130").ClearContents()
$ws.Range("I130").ClearContents()
$ws.Range("J130").ClearContents()
$ws.Range("K130").ClearContents()
$ws.Range("L130").ClearContents()
$ws.Range("H131").ClearContents()
$ws.Range("I131").ClearContents()
$ws.Range("J131").ClearContents()
$ws.Range("K131").ClearContents()
$ws.Range("L131").ClearContents()
$ws.Range("M131").ClearContents()
$ws.Range("N131").ClearContents()
$ws.Range("H132").ClearContents()
$ws.Range("I132").ClearContents()
$ws.Range("J132").ClearContents()
$ws.Range("K132").ClearContents()
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H133").ClearContents()
$ws.Range("I133").ClearContents()
$ws.Range("J133").ClearContents()
$ws.Range("K133").ClearContents()
$ws.Range("L133").ClearContents()
$ws.Range("H134").ClearContents()
$ws.Range("I134").ClearContents()
$ws.Range("J134").ClearContents()
$ws.Range("K134").ClearContents()
$ws.Range("L134").ClearContents()
$ws.Range("N134").ClearContents()
$ws.Range("H135").ClearContents()
$ws.Range("I135").ClearContents()
$ws.Range("J135").ClearContents()
$ws.Range("K135").ClearContents()
$ws.Range("L135").ClearContents()
$ws.Range("M135").ClearContents()
$ws.Range("H136").ClearContents()
$ws.Range("I136").ClearContents()
$ws.Range("J136").ClearContents()
$ws.Range("K136").ClearContents()
$ws.Range("L136").ClearContents()
$ws.Range("N136").ClearContents()
$ws.Range("H137").ClearContents()
$ws.Range("I137").ClearContents()
$ws.Range("J137").ClearContents()
$ws.Range("K137").ClearContents()
$ws.Range("L137").ClearContents()
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()
$ws.Range("H138").ClearContents()
$ws.Range("I138").ClearContents()
$ws.Range("J138").ClearContents()
$ws.Range("K138").ClearContents()
$ws.Range("L138").ClearContents()
$ws.Range("M138").ClearContents()
$ws.Range("N138").ClearContents()
$ws.Range("H139").ClearContents()
$ws.Range("I139").ClearContents()
$ws.Range("J139").ClearContents()
$ws.Range("K139").ClearContents()
$ws.Range("L139").ClearContents()
$ws.Range("N139").ClearContents()
$ws.Range("H140").ClearContents()
$ws.Range("I140").ClearContents()
$ws.Range("J140").ClearContents()
$ws.Range("K140").ClearContents()
$ws.Range("L140").ClearContents()
$ws.Range("N140").ClearContents()
$ws.Range("H141").ClearContents()
$ws.Range("I141").ClearContents()
$ws.Range("J141").ClearContents()
$ws.Range("K141").ClearContents()
$ws.Range("L141").ClearContents()
$ws.Range("M141").ClearContents()

# --- Sheet ARM (index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = [double]3676.5356
$ws.Range("I32").Value = [double]3997.682
$ws.Range("J32").Value = [double]2499
$ws.Range("K32").Value = [double]3997.682
$ws.Range("L32").Value = [double]2499
$ws.Range("M32").Value = [double]-3710.682
$ws.Range("N32").Value = [double]-3073
$ws.Range("H63").Value = [double]9745.210999999999
$ws.Range("I63").Value = [double]10461.5
$ws.Range("J63").Value = [double]7739.6
$ws.Range("K63").Value = [double]10461.5
$ws.Range("L63").Value = [double]7739.6
$ws.Range("M63").Value = [double]-9775.5
$ws.Range("N63").Value = [double]-9111.6
$ws.Range("H66").Value = [double]9745.210999999999
$ws.Range("I66").Value = [double]10461.5
$ws.Range("J66").Value = [double]7739.6
$ws.Range("K66").Value = [double]52307.5
$ws.Range("L66").Value = [double]38698
$ws.Range("M66").Value = [double]-48875.5
$ws.Range("N66").Value = [double]-45562
$ws.Range("H74").Value = [double]990.4
$ws.Range("I74").Value = [double]850.8182
$ws.Range("K74").Value = [double]850.8182
$ws.Range("M74").Value = [double]23.18179999999995
$ws.Range("H77").Value = [double]990.4
$ws.Range("I77").Value = [double]850.8182
$ws.Range("K77").Value = [double]4254.091
$ws.Range("M77").Value = [double]113.9089999999997
$ws.Range("H102").Value = [double]3236.75
$ws.Range("I102").Value = [double]2842
$ws.Range("K102").Value = [double]2842
$ws.Range("M102").Value = [double]-1220
$ws.Range("H122").Value = [double]2449.75
$ws.Range("J122").Value = [double]2933
$ws.Range("L122").Value = [double]8799
$ws.Range("N122").Value = [double]-13699
$ws.Range("H132").Value = [double]3853
$ws.Range("I132").Value = [double]3832.875
$ws.Range("K132").Value = [double]11498.625
$ws.Range("M132").Value = [double]-8968.625

# --- Sheet BSM (index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = [double]1502.2222
$ws.Range("I20").Value = [double]1517.2858
$ws.Range("J20").Value = [double]1449.5
$ws.Range("K20").Value = [double]1517.2858
$ws.Range("L20").Value = [double]1449.5
$ws.Range("M20").Value = [double]-1270.2858
$ws.Range("N20").Value = [double]-1943.5
$ws.Range("H86").Value = [double]5499.6665
$ws.Range("I86").Value = [double]5249.5
$ws.Range("J86").Value = [double]6000
$ws.Range("K86").Value = [double]5249.5
$ws.Range("L86").Value = [double]6000
$ws.Range("M86").Value = [double]-4126.5
$ws.Range("N86").Value = [double]-8246
$ws.Range("H89").Value = [double]5499.6665
$ws.Range("I89").Value = [double]5249.5
$ws.Range("J89").Value = [double]6000
$ws.Range("K89").Value = [double]26247.5
$ws.Range("L89").Value = [double]30000
$ws.Range("M89").Value = [double]-20631.5
$ws.Range("N89").Value = [double]-41232
$ws.Range("H94").Value = [double]1078.1
$ws.Range("I94").Value = [double]799
$ws.Range("J94").Value = [double]1496.75
$ws.Range("K94").Value = [double]799
$ws.Range("L94").Value = [double]1496.75
$ws.Range("M94").Value = [double]-348
$ws.Range("N94").Value = [double]-2398.75
$ws.Range("H99").Value = [double]0
$ws.Range("I99").Value = [double]0
$ws.Range("K99").Value = [double]0
$ws.Range("M99").ClearContents()
$ws.Range("H107").Value = [double]3744.1052
$ws.Range("I107").Value = [double]3794.4
$ws.Range("J107").Value = [double]3555.5
$ws.Range("K107").Value = [double]3794.4
$ws.Range("L107").Value = [double]3555.5
$ws.Range("M107").Value = [double]-1874.4
$ws.Range("N107").Value = [double]-7395.5
$ws.Range("H117").ClearContents()
$ws.Range("I117").ClearContents()
$ws.Range("J117").ClearContents()
$ws.Range("K117").ClearContents()
$ws.Range("L117").ClearContents()
$ws.Range("H118").ClearContents()
$ws.Range("I118").ClearContents()
$ws.Range("J118").ClearContents()
$ws.Range("K118").ClearContents()
$ws.Range("L118").ClearContents()
$ws.Range("H119").ClearContents()
$ws.Range("I119").ClearContents()
$ws.Range("J119").ClearContents()
$ws.Range("K119").ClearContents()
$ws.Range("L119").ClearContents()
$ws.Range("N119").ClearContents()
$ws.Range("H120").ClearContents()
$ws.Range("I120").ClearContents()
$ws.Range("J120").ClearContents()
$ws.Range("K120").ClearContents()
$ws.Range("L120").ClearContents()
$ws.Range("H122").ClearContents()
$ws.Range("I122").ClearContents()
$ws.Range("J122").ClearContents()
$ws.Range("K122").ClearContents()
$ws.Range("L122").ClearContents()
$ws.Range("H123").ClearContents()
$ws.Range("I123").ClearContents()
$ws.Range("J123").ClearContents()
$ws.Range("K123").ClearContents()
$ws.Range("L123").ClearContents()
$ws.Range("H124").ClearContents()
$ws.Range("I124").ClearContents()
$ws.Range("J124").ClearContents()
$ws.Range("K124").ClearContents()
$ws.Range("L124").ClearContents()
$ws.Range("H125").ClearContents()
$ws.Range("I125").ClearContents()
$ws.Range("J125").ClearContents()
$ws.Range("K125").ClearContents()
$ws.Range("L125").ClearContents()
$ws.Range("N125").ClearContents()
$ws.Range("H126").ClearContents()
$ws.Range("I126").ClearContents()
$ws.Range("J126").ClearContents()
$ws.Range("K126").ClearContents()
$ws.Range("L126").ClearContents()
$ws.Range("H127").ClearContents()
$ws.Range("I127").ClearContents()
$ws.Range("J127").ClearContents()
$ws.Range("K127").ClearContents()
$ws.Range("L127").ClearContents()
$ws.Range("M127").ClearContents()
$ws.Range("H128").ClearContents()
$ws.Range("I128").ClearContents()
$ws.Range("J128").ClearContents()
$ws.Range("K128").ClearContents()
$ws.Range("L128").ClearContents()
$ws.Range("H129").ClearContents()
$ws.Range("I129").ClearContents()
$ws.Range("J129").ClearContents()
$ws.Range("K129").ClearContents()
$ws.Range("L129").ClearContents()
$ws.Range("N129").ClearContents()
$ws.Range("H130").ClearContents()
$ws.Range("I130").ClearContents()
$ws.Range("J130").ClearContents()
$ws.Range("K130").ClearContents()
$ws.Range("L130").ClearContents()
$ws.Range("N130").ClearContents()
$ws.Range("H131").ClearContents()
$ws.Range("I131").ClearContents()
$ws.Range("J131").ClearContents()
$ws.Range("K131").ClearContents()
$ws.Range("L131").ClearContents()
$ws.Range("H132").ClearContents()
$ws.Range("I132").ClearContents()
$ws.Range("J132").ClearContents()
$ws.Range("K132").ClearContents()
$ws.Range("L132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H133").ClearContents()
$ws.Range("I133").ClearContents()
$ws.Range("J133").ClearContents()
$ws.Range("K133").ClearContents()
$ws.Range("L133").ClearContents()
$ws.Range("M133").ClearContents()
$ws.Range("N133").ClearContents()
$ws.Range("H134").ClearContents()
$ws.Range("I134").ClearContents()
$ws.Range("J134").ClearContents()
$ws.Range("K134").ClearContents()
$ws.Range("L134").ClearContents()
$ws.Range("M134").ClearContents()
$ws.Range("H135").ClearContents()
$ws.Range("I135").ClearContents()
$ws.Range("J135").ClearContents()
$ws.Range("K135").ClearContents()
$ws.Range("L135").ClearContents()
$ws.Range("N135").ClearContents()
$ws.Range("H137").ClearContents()
$ws.Range("I137").ClearContents()
$ws.Range("J137").ClearContents()
$ws.Range("K137").ClearContents()
$ws.Range("L137").ClearContents()
$ws.Range("H138").ClearContents()
$ws.Range("I138").ClearContents()
$ws.Range("J138").ClearContents()
$ws.Range("K138").ClearContents()
$ws.Range("L138").ClearContents()
$ws.Range("H139").ClearContents()
$ws.Range("I139").ClearContents()
$ws.Range("J139").ClearContents()
$ws.Range("K139").ClearContents()
$ws.Range("L139").ClearContents()
$ws.Range("H140").ClearContents()
$ws.Range("I140").ClearContents()
$ws.Range("J140").ClearContents()
$ws.Range("K140").ClearContents()
$ws.Range("L140").ClearContents()
$ws.Range("N140").ClearContents()
$ws.Range("H141").ClearContents()
$ws.Range("I141").ClearContents()
$ws.Range("J141").ClearContents()
$ws.Range("K141").ClearContents()
$ws.Range("L141").ClearContents()

# --- Sheet CRP (index 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value = [double]1000
$ws.Range("I22").Value = [double]1000
$ws.Range("J22").Value = [double]1000
$ws.Range("K22").Value = [double]1000
$ws.Range("L22").Value = [double]1000
$ws.Range("M22").Value = [double]-650
$ws.Range("N22").Value = [double]-1700
$ws.Range("H31").Value = [double]2417.05
$ws.Range("I31").Value = [double]2849.7273
$ws.Range("J31").Value = [double]1888.2222
$ws.Range("K31").Value = [double]2849.7273
$ws.Range("L31").Value = [double]1888.2222
$ws.Range("M31").Value = [double]-2554.7273
$ws.Range("N31").Value = [double]-2478.2222
$ws.Range("H34").Value = [double]2417.05
$ws.Range("I34").Value = [double]2849.7273
$ws.Range("J34").Value = [double]1888.2222
$ws.Range("K34").Value = [double]2849.7273
$ws.Range("L34").Value = [double]1888.2222
$ws.Range("M34").Value = [double]-2647.7273
$ws.Range("N34").Value = [double]-2292.2222
$ws.Range("H56").Value = [double]6531
$ws.Range("I56").Value = [double]6531
$ws.Range("J56").Value = [double]0
$ws.Range("K56").Value = [double]6531
$ws.Range("L56").Value = [double]0
$ws.Range("M56").ClearContents()
$ws.Range("N56").Value = [double]-5686
$ws.Range("H62").Value = [double]4148.9
$ws.Range("I62").Value = [double]3979.8
$ws.Range("J62").Value = [double]4318
$ws.Range("K62").Value = [double]3979.8
$ws.Range("L62").Value = [double]4318
$ws.Range("M62").Value = [double]-3355.8
$ws.Range("N62").Value = [double]-5566
$ws.Range("H65").Value = [double]4148.9
$ws.Range("I65").Value = [double]3979.8
$ws.Range("J65").Value = [double]4318
$ws.Range("K65").Value = [double]19899
$ws.Range("L65").Value = [double]21590
$ws.Range("M65").Value = [double]-16779
$ws.Range("N65").Value = [double]-27830
$ws.Range("H107").Value = [double]1223.6316
$ws.Range("I107").Value = [double]1144.1538
$ws.Range("K107").Value = [double]1144.1538
$ws.Range("M107").Value = [double]775.8462
$ws.Range("H129").ClearContents()
$ws.Range("I129").ClearContents()
$ws.Range("J129").ClearContents()
$ws.Range("K129").ClearContents()
$ws.Range("L129").ClearContents()
$ws.Range("H130").ClearContents()
$ws.Range("I130").ClearContents()
$ws.Range("J130").ClearContents()
$ws.Range("K130").ClearContents()
$ws.Range("L130").ClearContents()
$ws.Range("N130").ClearContents()
$ws.Range("H131").ClearContents()
$ws.Range("I131").ClearContents()
$ws.Range("J131").ClearContents()
$ws.Range("K131").ClearContents()
$ws.Range("L131").ClearContents()
$ws.Range("H132").ClearContents()
$ws.Range("I132").ClearContents()
$ws.Range("J132").ClearContents()
$ws.Range("K132").ClearContents()
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H133").ClearContents()
$ws.Range("I133").ClearContents()
$ws.Range("J133").ClearContents()
$ws.Range("K133").ClearContents()
$ws.Range("L133").ClearContents()
$ws.Range("M133").ClearContents()
$ws.Range("N133").ClearContents()
$ws.Range("H134").ClearContents()
$ws.Range("I134").ClearContents()
$ws.Range("J134").ClearContents()
$ws.Range("K134").ClearContents()
$ws.Range("L134").ClearContents()
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()
$ws.Range("H135").ClearContents()
$ws.Range("I135").ClearContents()
$ws.Range("J135").ClearContents()
$ws.Range("K135").ClearContents()
$ws.Range("L135").ClearContents()
$ws.Range("H137").ClearContents()
$ws.Range("I137").ClearContents()
$ws.Range("J137").ClearContents()
$ws.Range("K137").ClearContents()
$ws.Range("L137").ClearContents()
$ws.Range("H138").ClearContents()
$ws.Range("I138").ClearContents()
$ws.Range("J138").ClearContents()
$ws.Range("K138").ClearContents()
$ws.Range("L138").ClearContents()
$ws.Range("N138").ClearContents()
$ws.Range("H139").ClearContents()
$ws.Range("I139").ClearContents()
$ws.Range("J139").ClearContents()
$ws.Range("K139").ClearContents()
$ws.Range("L139").ClearContents()
$ws.Range("H140").ClearContents()
$ws.Range("I140").ClearContents()
$ws.Range("J140").ClearContents()
$ws.Range("K140").ClearContents()
$ws.Range("L140").ClearContents()
$ws.Range("N140").ClearContents()
$ws.Range("H141").ClearContents()
$ws.Range("I141").ClearContents()
$ws.Range("J141").ClearContents()
$ws.Range("K141").ClearContents()
$ws.Range("L141").ClearContents()
$ws.Range("M141").ClearContents()
$ws.Range("N141").ClearContents()

# --- Sheet CUL (index 5) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H60").Value = [double]1542
$ws.Range("I60").Value = [double]499
$ws.Range("K60").Value = [double]1497
$ws.Range("M60").Value = [double]-1246
$ws.Range("H68").Value = [double]1643.5
$ws.Range("I68").Value = [double]1000
$ws.Range("J68").Value = [double]1735.4286
$ws.Range("K68").Value = [double]3000
$ws.Range("L68").Value = [double]5206.2858
$ws.Range("M68").Value = [double]-2189
$ws.Range("N68").Value = [double]-6828.2858
$ws.Range("H71").Value = [double]1643.5
$ws.Range("I71").Value = [double]1000
$ws.Range("J71").Value = [double]1735.4286
$ws.Range("K71").Value = [double]9000
$ws.Range("L71").Value = [double]15618.8574
$ws.Range("M71").Value = [double]-4944
$ws.Range("N71").Value = [double]-23730.8574
$ws.Range("H102").Value = [double]5000
$ws.Range("I102").Value = [double]5000
$ws.Range("K102").Value = [double]15000
$ws.Range("M102").Value = [double]-12566
$ws.Range("H131").Value = [double]835775.75
$ws.Range("I131").Value = [double]2059
$ws.Range("K131").Value = [double]6177
$ws.Range("M131").Value = [double]-1137

# --- Sheet GSM (index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = [double]7199.6
$ws.Range("I70").Value = [double]6999.75
$ws.Range("K70").Value = [double]6999.75
$ws.Range("M70").Value = [double]-6729.75
$ws.Range("H73").Value = [double]7199.6
$ws.Range("I73").Value = [double]6999.75
$ws.Range("K73").Value = [double]6999.75
$ws.Range("M73").Value = [double]-6063.75
$ws.Range("H80").Value = [double]4499.5
$ws.Range("I80").Value = [double]4499.5
$ws.Range("K80").Value = [double]4499.5
$ws.Range("M80").Value = [double]-3501.5
$ws.Range("H83").Value = [double]4499.5
$ws.Range("I83").Value = [double]4499.5
$ws.Range("K83").Value = [double]22497.5
$ws.Range("M83").Value = [double]-17505.5
$ws.Range("H97").Value = [double]1162
$ws.Range("I97").Value = [double]1610.3334
$ws.Range("J97").Value = [double]489.5
$ws.Range("K97").Value = [double]1610.3334
$ws.Range("L97").Value = [double]489.5
$ws.Range("M97").Value = [double]-1114.3334
$ws.Range("N97").Value = [double]-1481.5
$ws.Range("H126").Value = [double]2499.6667
$ws.Range("I126").Value = [double]2499.6667
$ws.Range("K126").Value = [double]7499.000100000001
$ws.Range("M126").Value = [double]-5029.000100000001
$ws.Range("H132").Value = [double]4595.4
$ws.Range("I132").Value = [double]2992.6667
$ws.Range("J132").Value = [double]6999.5
$ws.Range("K132").Value = [double]8978.000100000001
$ws.Range("L132").Value = [double]20998.5
$ws.Range("M132").Value = [double]-6448.000100000001
$ws.Range("N132").Value = [double]-26058.5

# --- Sheet LTW (index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H22").Value = [double]1974.8334
$ws.Range("I22").Value = [double]1712.25
$ws.Range("K22").Value = [double]1712.25
$ws.Range("M22").Value = [double]-1417.25
$ws.Range("H27").Value = [double]1974.8334
$ws.Range("I27").Value = [double]1712.25
$ws.Range("K27").Value = [double]1712.25
$ws.Range("M27").Value = [double]-1605.25
$ws.Range("H40").Value = [double]3981.8333
$ws.Range("I40").Value = [double]3981.8333
$ws.Range("K40").Value = [double]3981.8333
$ws.Range("M40").Value = [double]-3845.8333
$ws.Range("H68").Value = [double]2915
$ws.Range("J68").Value = [double]2836
$ws.Range("L68").Value = [double]2836
$ws.Range("N68").Value = [double]-4334
$ws.Range("H71").Value = [double]2915
$ws.Range("J71").Value = [double]2836
$ws.Range("L71").Value = [double]14180
$ws.Range("N71").Value = [double]-21668
$ws.Range("H82").Value = [double]1081.1818
$ws.Range("J82").Value = [double]1492.25
$ws.Range("L82").Value = [double]1492.25
$ws.Range("N82").Value = [double]-2214.25
$ws.Range("H85").Value = [double]1081.1818
$ws.Range("J85").Value = [double]1492.25
$ws.Range("L85").Value = [double]1492.25
$ws.Range("N85").Value = [double]-3988.25
$ws.Range("H93").Value = [double]3800
$ws.Range("I93").Value = [double]0
$ws.Range("J93").Value = [double]3800
$ws.Range("K93").Value = [double]0
$ws.Range("L93").ClearContents()
$ws.Range("M93").Value = [double]3800
$ws.Range("N93").Value = [double]-6296
$ws.Range("H100").Value = [double]1696.0625
$ws.Range("I100").Value = [double]1722.4667
$ws.Range("J100").Value = [double]1300
$ws.Range("K100").Value = [double]1722.4667
$ws.Range("L100").Value = [double]1300
$ws.Range("M100").Value = [double]-1181.4667
$ws.Range("N100").Value = [double]-2382
$ws.Range("H132").Value = [double]5554.25
$ws.Range("I132").Value = [double]5554.25
$ws.Range("K132").Value = [double]16662.75
$ws.Range("M132").Value = [double]-14132.75

# --- Sheet WVR (index 8) ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H81").Value = [double]909742.75
$ws.Range("I81").Value = [double]696.44446
$ws.Range("K81").Value = [double]1392.88892
$ws.Range("M81").Value = [double]-331.8889200000001
$ws.Range("H84").Value = [double]909742.75
$ws.Range("I84").Value = [double]696.44446
$ws.Range("K84").Value = [double]6964.444600000001
$ws.Range("M84").Value = [double]-1660.444600000001
$ws.Range("H96").Value = [double]999
$ws.Range("J96").Value = [double]0
$ws.Range("L96").Value = [double]0
$ws.Range("N96").ClearContents()
$ws.Range("H132").Value = [double]5060.8423
$ws.Range("I132").Value = [double]5060.8423
$ws.Range("K132").Value = [double]15182.5269
$ws.Range("M132").Value = [double]-12652.5269
